$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.888.30'
$ws.Range("E2").Value = '  -2.27%  '
$ws.Range("D3").Value = '1.652.91'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3893'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("E8").Value = '  -2.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.345'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08466'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.047'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.057'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001315'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("D17").Value = '1.654.65'
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06986'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("E20").Value = '  -4.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.992'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").Value = '23.897.67'
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.430'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.973'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.424'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '137.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.30%  '
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.483'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.56%  '
$ws.Range("D33").Value = '1.840.64'
$ws.Range("E33").Value = '  -0.33%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.006'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.26%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08130'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.73%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.686'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.27%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02910'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2675'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09122'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7564'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.425'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.55'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("E45").Value = '  -2.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.449'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.112'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9997'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08278'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.232'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.08%  '
